# Generate Report for Handoff
#
# The previous handoff round (source id 26e25ad3-b904-4ce2-b74f-15b9e6e917d7)
# has been superseded by a new handoff round (source id
# 9dd30674-4adc-447b-b498-6616376d2697, new target-xlf checksum
# 99803a9f9145700dae47ef704869ab13733de735). Update the status report so
# every sheet reflects the new source/target file names and the new
# handoff timestamps recorded for this round.

$wb = $excel.ActiveWorkbook

$oldId = "26e25ad3-b904-4ce2-b74f-15b9e6e917d7"
$newId = "9dd30674-4adc-447b-b498-6616376d2697"
$oldHash = "4e5441ad70ba7a6633c9361a3b45aa30aad2ae66"
$newHash = "99803a9f9145700dae47ef704869ab13733de735"

# --- Overview sheet: just the handed-off markdown source file name ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"

# --- zh-cn sheet: source file, latest handoff (xlf) file, handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("C2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-03-09 15:51:19"

# --- de-de sheet: source file, latest handoff (xlf) file, handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("C2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-03-09 15:51:28"
